# Test Matrix.xlsx edit script
# - "crcN fully functional, removed unneeded input and cleaned up comments"
# - Highlights the working CRC results (green fill) on the "Short Form Testing"
#   sheet, adds a right-aligned+filled style to the winning-polynomial row,
#   clears the now-unneeded 16-bit example row, and updates the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Short Form Testing")

# Green fill already present in the workbook as fillId=3 (RGB 92D050).
# Using this exact color lets Excel reuse/create the matching cellXfs entry.
$green = 5296274   # BGR for RGB(146, 208, 80) / hex 92D050

# Rows 3-12: columns F:L (incl. K) get the green "valid result" fill.
$ws.Range("F3:L12").Interior.Color = $green

# Row 13 is the "winner" row (already right aligned); keep the alignment but
# add the same green fill across F:J (K stays unfilled, L already has its
# own grey fill from before).
$ws.Range("F13:J13").Interior.Color = $green

# Rows 14-15: columns F,G,H,I,J get the green fill (K stays unfilled).
$ws.Range("F14:J15").Interior.Color = $green

# Row 16 (16-bit example) is no longer needed - blank out its results,
# keeping the existing (unfilled) cell style.
$ws.Range("F16:K16").ClearContents()

# Restore the saved scroll position / selection for the sheet.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("J12").Select()
